# The commit changes a single data cell in the "Rules" decision table:
# cell C10 (the lower bound "min" for rule R20) is updated from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
